$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.258.96"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.55%  "
$ws.Range("D3").Value = "'2.421.02"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.94%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'562.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.69%  "
$ws.Range("D6").Value = "'143.98"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.25%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +1.49%  "
$ws.Range("D9").Value = "'2.418.92"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.79%  "
$ws.Range("E10").Value = "  +1.78%  "
$ws.Range("E11").Value = "  -2.13%  "
$ws.Range("E12").Value = "  +0.16%  "
$ws.Range("E13").Value = "  +0.76%  "
$ws.Range("D14").Value = "'25.88"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.69%  "
$ws.Range("E15").Value = "  +4.11%  "
$ws.Range("D16").Value = "'2.859.14"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.00%  "
$ws.Range("D17").Value = "'62.162.20"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.50%  "
$ws.Range("D18").Value = "'2.419.63"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.99%  "
$ws.Range("E19").Value = "  +2.93%  "
$ws.Range("E20").Value = "  +1.10%  "
$ws.Range("D21").Value = "'323.92"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.20%  "
$ws.Range("E22").Value = "  +0.80%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").Value = "'65.58"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.17%  "
$ws.Range("E25").Value = "  -2.75%  "
$ws.Range("D26").Value = "'8.94"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.99%  "
$ws.Range("D27").Value = "'582.35"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +12.72%  "
$ws.Range("D28").Value = "'2.541.19"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.97%  "
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("D30").Value = "'0.0₃0940"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.60%  "
$ws.Range("E31").Value = "  +4.52%  "
$ws.Range("D32").Value = "'8.24"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.02%  "
$ws.Range("D33").Value = "'0.151"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.93%  "
$ws.Range("E34").Value = "  +2.35%  "
$ws.Range("E35").Value = "  +1.97%  "
$ws.Range("D36").Value = "'5.70"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.55%  "
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("E38").Value = "  +2.23%  "
$ws.Range("E39").Value = "  +1.64%  "
$ws.Range("D40").Value = "'152.53"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.72%  "
$ws.Range("D41").Value = "'18.64"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.95%  "
$ws.Range("E42").Value = "  -3.57%  "
$ws.Range("E43").Value = "  -0.17%  "
$ws.Range("E44").Value = "  +8.50%  "
$ws.Range("D45").Value = "'150.11"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.92%  "
$ws.Range("D46").Value = "'3.65"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.67%  "
$ws.Range("D47").Value = "'0.0537"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.73%  "
$ws.Range("D48").Value = "'20.30"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.26%  "
$ws.Range("D50").Value = "'0.0924"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.92%  "
$ws.Range("E51").Value = "  +2.09%  "
